$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# --- Shape 52: "username" ellipse (cNvPr id=161) -> underline the text ---
$shp = $s.Shapes.Item(52)
$shp.TextFrame.TextRange.Font.Underline = -1

# --- Shape 57: "id" ellipse (cNvPr id=60) ---
# move up slightly (y: 40042 -> 13409 EMU) and underline the text
$shp = $s.Shapes.Item(57)
$shp.Top = 1.0558267831802368
$shp.TextFrame.TextRange.Font.Underline = -1

# --- Shape 58: connector line (cNvPr id=7) linking shape 57 to shape below ---
# off y: 502754 -> 476121 EMU ; ext cy: 287476 -> 314109 EMU
$shp = $s.Shapes.Item(58)
$shp.Top = 37.489845275878906
$shp.Height = 24.73299217224121

# --- Shape 59: "id" ellipse (cNvPr id=67) -> underline the text ---
$shp = $s.Shapes.Item(59)
$shp.TextFrame.TextRange.Font.Underline = -1

# --- Shape 61: "id" ellipse (cNvPr id=72) -> underline the text ---
$shp = $s.Shapes.Item(61)
$shp.TextFrame.TextRange.Font.Underline = -1

# --- Shape 63: "id" ellipse (cNvPr id=74) -> underline the text ---
$shp = $s.Shapes.Item(63)
$shp.TextFrame.TextRange.Font.Underline = -1
